# Applies the "art of bop drumming" edit:
#  - Renames Hoja1 -> agend, Hoja2 -> dias
#  - Makes "dias" (Hoja2) the active/selected sheet (activeTab=1, tabSelected on sheet2)
#  - On "agend" (Hoja1): fills G2 with a new note, appends ", " to G3's note,
#    and moves the sheet's saved selection to G2 (no longer the tab-selected sheet)
#  - On "dias" (Hoja2): builds a small weekly schedule grid

$wb = $excel.ActiveWorkbook

# --- Rename sheets -------------------------------------------------------
$wsAgend = $wb.Worksheets.Item("Hoja1")
$wsAgend.Name = "agend"

$wsDias = $wb.Worksheets.Item("Hoja2")
$wsDias.Name = "dias"

# --- agend (sheet1) content tweaks ---------------------------------------
$wsAgend.Range("G3").Value = "diogenes, helenismo, "
$wsAgend.Range("G2").Value = "teogonía, homero/McIntyre?"

# Update the saved selection on agend to G2 and make it no longer the
# active/tab-selected sheet.
$wsAgend.Range("G2").Select()

# --- dias (sheet2) content ------------------------------------------------
$wsDias.Range("B2").Value = "lunes"
$wsDias.Range("C2").Value = "martes"
$wsDias.Range("D2").Value = "miércoles"
$wsDias.Range("E2").Value = "jueves"
$wsDias.Range("F2").Value = "viernes"
$wsDias.Range("G2").Value = "sábado"

$wsDias.Range("G4").Value = "P. I. H."
$wsDias.Range("A4").Value = "7:20 - 9:00"

$wsDias.Range("A10").Value = "13:00 - 14:20"

$wsDias.Range("F4").Value = "PENS C"
$wsDias.Range("F10").Value = "PENS C"

# Make dias the active sheet/tab and set its saved selection.
$wsDias.Activate()
$wsDias.Range("F11").Select()
